{"js": "// Replace the 25 two-digit-division answer cells in the single table with\n// a new set of problems, preserving all cell/run/paragraph formatting.\n// The grid is read via Table.values (row-major) and written back the same\n// way so only the <w:t> text content changes.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Old -> new text, keyed by (row, col) position within the table grid.\n// Only the five \"data\" rows (0, 4, 8, 12, 16) carry visible text; the\n// other rows are intentionally blank spacer rows and are left untouched.\nconst newGridByRow = {\n  0: [\"20\u00f73=6, 2\", \"25\u00f72=12, 1\", \"95\u00f79=10, 5\", \"93\u00f77=13, 2\", \"80\u00f76=13, 2\"],\n  4: [\"12\u00f78=1, 4\", \"39\u00f79=4, 3\", \"22\u00f79=2, 4\", \"61\u00f74=15, 1\", \"16\u00f74=4, 0\"],\n  8: [\"37\u00f74=9, 1\", \"97\u00f73=32, 1\", \"43\u00f73=14, 1\", \"64\u00f72=32, 0\", \"90\u00f74=22, 2\"],\n  12: [\"78\u00f77=11, 1\", \"22\u00f75=4, 2\", \"93\u00f73=31, 0\", \"52\u00f73=17, 1\", \"48\u00f75=9, 3\"],\n  16: [\"97\u00f79=10, 7\", \"83\u00f78=10, 3\", \"79\u00f77=11, 2\", \"12\u00f75=2, 2\", \"87\u00f78=10, 7\"],\n};\n\nconst values = table.values.map((row) => row.slice());\n\nfor (const rowIndexStr of Object.keys(newGridByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const newRow = newGridByRow[rowIndex];\n  for (let col = 0; col < newRow.length; col++) {\n    values[rowIndex][col] = newRow[col];\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-division answer cells in the single table with\n# a new set of problems, preserving all cell/run/paragraph formatting.\n# Only the five \"data\" rows (1, 5, 9, 13, 17 in 1-based COM indexing)\n# carry visible text; the other rows are blank spacer rows and are left\n# untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"20\u00f73=6, 2\", \"25\u00f72=12, 1\", \"95\u00f79=10, 5\", \"93\u00f77=13, 2\", \"80\u00f76=13, 2\")\n    5  = @(\"12\u00f78=1, 4\", \"39\u00f79=4, 3\", \"22\u00f79=2, 4\", \"61\u00f74=15, 1\", \"16\u00f74=4, 0\")\n    9  = @(\"37\u00f74=9, 1\", \"97\u00f73=32, 1\", \"43\u00f73=14, 1\", \"64\u00f72=32, 0\", \"90\u00f74=22, 2\")\n    13 = @(\"78\u00f77=11, 1\", \"22\u00f75=4, 2\", \"93\u00f73=31, 0\", \"52\u00f73=17, 1\", \"48\u00f75=9, 3\")\n    17 = @(\"97\u00f79=10, 7\", \"83\u00f78=10, 3\", \"79\u00f77=11, 2\", \"12\u00f75=2, 2\", \"87\u00f78=10, 7\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Count; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
